# Remove duplicated teacher entries ("Listas sem duplicação de professores")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose duplicated class-list text should be cleared to "-"
$clearCells = @(
    "E3",
    "B4", "E4",
    "B6", "E6",
    "E7",
    "C11",
    "C12",
    "C14",
    "C15",
    "B18", "C18",
    "B19", "C19",
    "B20", "C20",
    "B21", "C21"
)

foreach ($cell in $clearCells) {
    $ws.Range($cell).Value = "-"
}

# The "MCT-1A-Tecnologia dos Materiais." list moves into F7 and F8
$ws.Range("F7").Value = "['MCT-1A-Tecnologia dos Materiais.', -]"
$ws.Range("F8").Value = "['MCT-1A-Tecnologia dos Materiais.', -]"
